# Apply marksheet corrections: update correct/total marks values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# B11: Marking row, Right column -> 5
$ws.Range("B11").Value = 5

# B12: Total row, Right column -> 80
$ws.Range("B12").Value = 80

# E12: Total row, Max column -> "80/140"
$ws.Range("E12").Value = "80/140"
